# Refresh the cryptos price table (rows 2-51) to the latest scrape.
# A new coin (WrappedeETH) was inserted at row 28, shifting all
# subsequent rows down by one; the last prior row (WhiteBITCoin)
# drops off the bottom of the fixed 50-row table as a result.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '60.824.98'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +0.14%  '

# Row 3: Ethereum
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.349.08'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -1.12%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  -0.03%  '

# Row 5: BNB
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '544.29'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.10%  '

# Row 6: Solana
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '136.66'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -3.14%  '

# Row 7: USDC
$ws.Range('E7').Value = '  +0.01%  '

# Row 8: XRP
$ws.Range('E8').Value = '  -1.99%  '

# Row 9: LidoStakedEther
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.347.66'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -1.18%  '

# Row 10: Dogecoin
$ws.Range('E10').Value = '  -0.16%  '

# Row 11: TRON
$ws.Range('E11').Value = '  +2.05%  '

# Row 12: Toncoin
$ws.Range('E12').Value = '  -0.43%  '

# Row 13: Cardano
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.343'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -0.08%  '

# Row 14: Avalanche
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '24.64'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -3.24%  '

# Row 15: WrappedliquidstakedEther2.0
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '2.771.53'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -1.16%  '

# Row 16: WrappedBTC
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '60.990.62'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.58%  '

# Row 17: ShibaInu
$ws.Range('E17').Value = '  -1.31%  '

# Row 18: WrappedEther
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.348.17'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -1.11%  '

# Row 19: Chainlink
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '10.62'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.13%  '

# Row 20: BitcoinCash
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '319.31'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +1.13%  '

# Row 21: Polkadot
$ws.Range('E21').Value = '  +0.57%  '

# Row 22: Uniswap
$ws.Range('E22').Value = '  -2.08%  '

# Row 23: Dai
$ws.Range('E23').Value = '  +0.02%  '

# Row 24: Litecoin
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '63.35'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +1.02%  '

# Row 25: SuiNetwork
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.67'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -8.02%  '

# Row 26: Aptos
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '8.40'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +9.07%  '

# Row 27: Binance-PegBSC-USD
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.06%  '

# Row 28: WrappedeETH
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.464.72'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -1.07%  '

# Row 29: InternetComputer(DFINITY)
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.02'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +0.70%  '

# Row 30: Fetch.AI
$ws.Range('E30').Value = '  -3.38%  '

# Row 31: Bittensor
$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '496.64'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -4.65%  '

# Row 32: PEPE
$ws.Range('B32').Value = 'PEPE'
$ws.Range('C32').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0₃0860'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -7.08%  '

# Row 33: Kaspa
$ws.Range('B33').Value = 'Kaspa'
$ws.Range('C33').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.146'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.85%  '

# Row 34: PancakeSwap
$ws.Range('B34').Value = 'PancakeSwap'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.79'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -2.31%  '

# Row 35: ImmutableX
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.51'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -3.04%  '

# Row 36: FirstDigitalUSD
$ws.Range('B36').Value = 'FirstDigitalUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.999'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.04%  '

# Row 37: NEARProtocol
$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '4.63'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.40%  '

# Row 38: PolygonEcosystemToken
$ws.Range('B38').Value = 'PolygonEcosystemToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.376'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +0.34%  '

# Row 39: EthereumClassic
$ws.Range('B39').Value = 'EthereumClassic'
$ws.Range('C39').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '18.50'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +3.05%  '

# Row 40: Stacks
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.82'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +5.74%  '

# Row 41: RenderToken
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.23'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -4.28%  '

# Row 42: Monero
$ws.Range('B42').Value = 'Monero'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '143.26'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +4.67%  '

# Row 43: USDe
$ws.Range('B43').Value = 'USDe'
$ws.Range('C43').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.999'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.09%  '

# Row 44: Aave
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '142.52'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +2.19%  '

# Row 45: Filecoin
$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.56'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.40%  '

# Row 46: dogwifhat
$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.02'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -8.90%  '

# Row 47: Hedera
$ws.Range('B47').Value = 'Hedera'
$ws.Range('C47').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0515'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.12%  '

# Row 48: InjectiveProtocol
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '19.12'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -5.66%  '

# Row 49: Mantle
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.567'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -1.23%  '

# Row 50: Stellar
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0900'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -1.20%  '

# Row 51: VeChain
$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0221'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.60%  '
